$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "37.449.86"
$c.ClearFormats()

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  +0.19%  "
$c.ClearFormats()

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.071.97"
$c.ClearFormats()

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  +0.58%  "
$c.ClearFormats()

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.00%  "
$c.ClearFormats()

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "235.23"
$c.ClearFormats()

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  -0.85%  "
$c.ClearFormats()

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.627"
$c.ClearFormats()

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +1.76%  "
$c.ClearFormats()

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  +0.02%  "
$c.ClearFormats()

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "57.32"
$c.ClearFormats()

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  -1.51%  "
$c.ClearFormats()

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.396"
$c.ClearFormats()

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  +3.20%  "
$c.ClearFormats()

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0773"
$c.ClearFormats()

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  +1.68%  "
$c.ClearFormats()

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  +0.89%  "
$c.ClearFormats()

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "2.378.12"
$c.ClearFormats()

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  +0.44%  "
$c.ClearFormats()

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "14.45"
$c.ClearFormats()

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +0.55%  "
$c.ClearFormats()

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "20.56"
$c.ClearFormats()

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -3.43%  "
$c.ClearFormats()

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.779"
$c.ClearFormats()

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  -0.01%  "
$c.ClearFormats()

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  -0.31%  "
$c.ClearFormats()

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.071.87"
$c.ClearFormats()

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  -1.29%  "
$c.ClearFormats()

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "37.381.87"
$c.ClearFormats()

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  -0.46%  "
$c.ClearFormats()

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.40"
$c.ClearFormats()

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  +3.98%  "
$c.ClearFormats()

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "69.82"
$c.ClearFormats()

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +1.24%  "
$c.ClearFormats()

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.0₃0819"
$c.ClearFormats()

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  +0.70%  "
$c.ClearFormats()

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "227.10"
$c.ClearFormats()

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  +0.80%  "
$c.ClearFormats()

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  +0.01%  "
$c.ClearFormats()

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  +1.04%  "
$c.ClearFormats()

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  -1.99%  "
$c.ClearFormats()

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "166.65"
$c.ClearFormats()

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +1.78%  "
$c.ClearFormats()

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "8.84"
$c.ClearFormats()

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -0.43%  "
$c.ClearFormats()

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.44"
$c.ClearFormats()

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  -0.60%  "
$c.ClearFormats()

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  +0.98%  "
$c.ClearFormats()

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "19.13"
$c.ClearFormats()

$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -0.45%  "
$c.ClearFormats()

$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -0.77%  "
$c.ClearFormats()

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "4.54"
$c.ClearFormats()

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  +0.98%  "
$c.ClearFormats()

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0620"
$c.ClearFormats()

$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  -1.22%  "
$c.ClearFormats()

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "4.58"
$c.ClearFormats()

$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  +2.21%  "
$c.ClearFormats()

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.49"
$c.ClearFormats()

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  -3.61%  "
$c.ClearFormats()

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.79"
$c.ClearFormats()

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +0.42%  "
$c.ClearFormats()

$c = $ws.Range("B37")
$c.NumberFormat = "@"
$c.Value = "RenderToken"
$c.ClearFormats()

$c = $ws.Range("C37")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c.ClearFormats()

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.31"
$c.ClearFormats()

$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  -1.38%  "
$c.ClearFormats()

$c = $ws.Range("B38")
$c.NumberFormat = "@"
$c.Value = "BinanceUSD"
$c.ClearFormats()

$c = $ws.Range("C38")
$c.NumberFormat = "@"
$c.Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$c.ClearFormats()

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.ClearFormats()

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  +0.01%  "
$c.ClearFormats()

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.73"
$c.ClearFormats()

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -2.06%  "
$c.ClearFormats()

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  -0.60%  "
$c.ClearFormats()

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "99.05"
$c.ClearFormats()

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  +1.35%  "
$c.ClearFormats()

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.0966"
$c.ClearFormats()

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  -1.78%  "
$c.ClearFormats()

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.486.89"
$c.ClearFormats()

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  +0.23%  "
$c.ClearFormats()

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.18"
$c.ClearFormats()

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  +2.63%  "
$c.ClearFormats()

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  +1.22%  "
$c.ClearFormats()

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "4.09"
$c.ClearFormats()

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  -9.84%  "
$c.ClearFormats()

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  +0.49%  "
$c.ClearFormats()

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "15.35"
$c.ClearFormats()

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  -3.75%  "
$c.ClearFormats()

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "7.26"
$c.ClearFormats()

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  +0.60%  "
$c.ClearFormats()

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  +0.86%  "
$c.ClearFormats()

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.264.81"
$c.ClearFormats()

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  +0.44%  "
$c.ClearFormats()
